$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10 (AIS): CurrentQ Q4 -> 3Q2024, Result Positive -> Negative
$ws.Range("C10").Value = "3Q2024"
$ws.Range("E10").Value = "Negative"

# Row 16 (Telstra): CurrentQ FY24 Half Year -> FY25 Half Year
$ws.Range("C16").Value = "FY25 Half Year"

# Row 20 (TRAI): Result Negative -> FetchingError
$ws.Range("E20").Value = "FetchingError"
